# Generate Report for Archive
# - Status moves from "Ready for handoff" to "In Translation" on every sheet
#   that tracks it (Overview!E/F and the zh-cn / de-de Status column), and
#   the now-shorter text lets those status columns narrow a bit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E:F").ColumnWidth = 12.43

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C:C").ColumnWidth = 12.43

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C:C").ColumnWidth = 12.43
